$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.957.42"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.68%  '

$ws.Range('D3').Value = "'1.906.38"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.99%  '

$ws.Range('D4').Value = "'0.9993"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.18%  '

$ws.Range('D5').Value = "'0.8073"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.20%  '

$ws.Range('D6').Value = "'241.39"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.95%  '

$ws.Range('E7').Value = '  -0.14%  '

$ws.Range('E8').Value = '  +2.88%  '

$ws.Range('D9').Value = "'26.44"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +4.82%  '

$ws.Range('D10').Value = "'0.07015"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.39%  '

$ws.Range('D11').Value = "'0.07999"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.82%  '

$ws.Range('D12').Value = "'1.909.43"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.87%  '

$ws.Range('D13').Value = "'0.7448"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.60%  '

$ws.Range('D14').Value = "'5.187"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.96%  '

$ws.Range('D15').Value = "'92.51"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.53%  '

$ws.Range('D16').Value = "'29.953.34"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.64%  '

$ws.Range('D17').Value = "'13.99"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.11%  '

$ws.Range('D18').Value = "'5.873"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.38%  '

$ws.Range('D19').Value = "'245.26"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.43%  '

$ws.Range('D20').Value = "'0.000007788"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.96%  '

$ws.Range('D21').Value = "'1.003"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.27%  '

$ws.Range('D22').Value = "'2.152.51"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.21%  '

$ws.Range('D23').Value = "'0.9995"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.21%  '

$ws.Range('D24').Value = "'6.938"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.88%  '

$ws.Range('D25').Value = "'168.12"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.69%  '

$ws.Range('D26').Value = "'9.216"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.39%  '

$ws.Range('D27').Value = "'0.1485"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +16.93%  '

$ws.Range('D28').Value = "'18.87"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.47%  '

$ws.Range('D29').Value = "'2.067"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.19%  '

$ws.Range('D30').Value = "'1.362"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.48%  '

$ws.Range('D31').Value = "'1.510"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.09%  '

$ws.Range('D32').Value = "'4.295"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.48%  '

$ws.Range('D33').Value = "'0.05534"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.17%  '

$ws.Range('D34').Value = "'4.070"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.24%  '

$ws.Range('D35').Value = "'1.261"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.51%  '

$ws.Range('D36').Value = "'0.7298"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.79%  '

$ws.Range('E37').Value = '  +0.17%  '

$ws.Range('D38').Value = "'0.01918"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.08%  '

$ws.Range('D39').Value = "'2.784"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.71%  '

$ws.Range('D40').Value = "'0.4405"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.83%  '

$ws.Range('D41').Value = "'72.09"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.31%  '

$ws.Range('D42').Value = "'5.972"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.44%  '

$ws.Range('D43').Value = "'0.9999"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.08%  '

$ws.Range('D44').Value = "'0.8387"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.77%  '

$ws.Range('D45').Value = "'1.889"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.77%  '

$ws.Range('D46').Value = "'100.90"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.66%  '

$ws.Range('D47').Value = "'7.569"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.37%  '

$ws.Range('D48').Value = "'9.707"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.32%  '

$ws.Range('D49').Value = "'982.58"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +9.62%  '

$ws.Range('D50').Value = "'2.058.60"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.42%  '

$ws.Range('D51').Value = "'36.16"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.91%  '
